$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E6 signature block text to the new multi-line signature block
$ws.Range("E6").Value = "Name: Alexandru Blagan .`n Title: Manager .`n Date:   "

# Re-fit the row height so it doesn't pick up an explicit custom row height
# (matches original workbook which had no explicit row height on row 6)
$ws.Rows.Item(6).AutoFit()

# Fill in rows 11-15 (A:E) with empty formatted cells that have the same
# border/style formatting as rows 7-10 (style index 7), without disturbing
# the vertical merges spanning rows 6-16 in columns A, B, C, D, E.
for ($r = 11; $r -le 15; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Borders.LineStyle = 1
    }
}
